$d = $word.ActiveDocument

# 1) Remove the existing "_GoBack" bookmark that sits after
#    "Remove email password from file"
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) Append a new list paragraph after the last paragraph
#    ("Optional but important: add captcha to register page")
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertAfter("`r" + "If you don" + [char]0x2019 + "t check the captcha for the first time, error msg appears. Then if you check and make a error in somewhere else, a new error msg and the captcha error msg appear.")

# 3) Re-insert a "_GoBack" bookmark as a zero-length point right before
#    the word "appear." in the newly added paragraph, splitting the run
#    the same way the original diff shows.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("appear.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkPoint = $d.Range($searchRange.Start, $searchRange.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint)
